$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-23 of column B use the date/time number format already present
# in the workbook (style index 1 applied to B2:B4 originally).
$dateFormat = "yyyy-mm-dd h:mm:ss"

# id, release_date (Excel serial date), tassativita
$rows = @(
    @(253295, 45911.58333333334, "X"),
    @(253549, 45912.58333333334, 0),
    @(253244, 45912.58333333334, 0),
    @(253377, 45910.58333333334, 0),
    @(252274, 45911.58333333334, 0),
    @(253371, 45910.58333333334, 0),
    @(253367, 45911.58333333334, 0),
    @(253392, 45911.58333333334, 0),
    @(253525, 45912.58333333334, 0),
    @(253376, 45910.58333333334, 0),
    @(253393, 45912.58333333334, 0),
    @(253472, 45912.58333333334, 0),
    @(253522, 45910.58333333334, 0),
    @(252529, 45910.58333333334, 0),
    @(252397, 45911.58333333334, 0),
    @(253528, 45911.58333333334, 0),
    @(253455, 45910.58333333334, 0),
    @(253436, 45910.58333333334, 0),
    @(253668, 45910.58333333334, 0),
    @(244743, 45910.58333333334, 0),
    @(253261, 45910.58333333334, 0),
    @(252741, 45911.58333333334, 0)
)

$rowIndex = 2
foreach ($r in $rows) {
    $ws.Range("A$rowIndex").Value2 = $r[0]
    $ws.Range("B$rowIndex").NumberFormat = $dateFormat
    $ws.Range("B$rowIndex").Value2 = $r[1]
    $ws.Range("C$rowIndex").Value = $r[2]
    $rowIndex++
}
